$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rename the sheet ----
$ws.Name = "Sheet1"

# ---- Header row (row 1) ----
$headers = @(
  "Date",
  "Model Name",
  "Exact Precision (Micro Avg)",
  "Exact Recall (Micro Avg)",
  "Exact F1 Score (Micro Avg)",
  "Exact Precision (Macro Avg)",
  "Exact Recall (Macro Avg)",
  "Exact F1 Score (Macro Avg)",
  "Exact Precision (Weighted Avg)",
  "Exact Recall (Weighted Avg)",
  "Exact F1 Score (Weighted Avg)",
  "Partial Precision",
  "Partial Recall",
  "Partial F1 Score",
  "Partial TP",
  "Partial FP",
  "Partial FN",
  "Support",
  "Accuracy",
  "Result Link",
  "Stats Link",
  "No of GPU Used",
  "Power Consumption"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Build the header format (bold font, centered + top aligned, thin box
# border all round) on a scratch cell far away from the used range, then
# paste just that *format* onto the header range in a single operation so
# only one extra style entry is produced for the whole header row.
$templateCell = $ws.Range("ZZ1000")
$templateCell.Value = "x"
$templateCell.Font.Bold = $true
$templateCell.HorizontalAlignment = -4108   # xlCenter
$templateCell.VerticalAlignment = -4160     # xlTop
$templateCell.Borders.LineStyle = 1         # xlContinuous

$headerRange = $ws.Range("A1:W1")
$templateCell.Copy()
$headerRange.PasteSpecial(-4122)            # xlPasteFormats
$templateCell.Clear()

# ---- Data row (row 2) ----
# Column A holds a date-formatted-looking string ("09/11/2025") that must
# stay literal text rather than being auto-converted into a date serial
# number. Writing it as a formula that evaluates to the text and then
# collapsing the formula down to its static result (copy / paste-values
# on the same cell) keeps it a plain string with no extra number format.
$dateCell = $ws.Cells.Item(2, 1)
$dateCell.Formula = '="09/11/2025"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)               # xlPasteValues

$ws.Cells.Item(2, 2).Value = "Qwen2.5-7B-Instruct"
$ws.Cells.Item(2, 3).Value = 0.448051948051948
$ws.Cells.Item(2, 4).Value = 0.2482014388489209
$ws.Cells.Item(2, 5).Value = 0.3194444444444445
$ws.Cells.Item(2, 6).Value = 0.2479724567959863
$ws.Cells.Item(2, 7).Value = 0.1226219961668878
$ws.Cells.Item(2, 8).Value = 0.1529871150845629
$ws.Cells.Item(2, 9).Value = 0.5086843547275202
$ws.Cells.Item(2, 10).Value = 0.2482014388489209
$ws.Cells.Item(2, 11).Value = 0.3179549041189638
$ws.Cells.Item(2, 12).Value = 0.5921052631578947
$ws.Cells.Item(2, 13).Value = 0.3249097472924188
$ws.Cells.Item(2, 14).Value = 0.4195804195804196
$ws.Cells.Item(2, 15).Value = 90
$ws.Cells.Item(2, 16).Value = 62
$ws.Cells.Item(2, 17).Value = 187
$ws.Cells.Item(2, 18).Value = 278
$ws.Cells.Item(2, 19).Value = 0.9540717869548437
$ws.Cells.Item(2, 20).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-7B-Instruct_3_shot.txt"
$ws.Cells.Item(2, 21).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-7B-Instruct_3_shot.txt"
$ws.Cells.Item(2, 22).Value = "4 MLGPU"
$ws.Cells.Item(2, 23).Value = "0.019 kWh"
$ws.Cells.Item(2, 24).Value = 674
